# Apply the nickname-table corrections described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text corrections -------------------------------------------------
$ws.Range("A2").Value  = "青行归"
$ws.Range("A18").Value = "PhuDgue"
$ws.Range("A35").Value = "Takion Kroslin"
$ws.Range("B15").Value = "何奕嘉"

# --- Avatar URL replacement (becomes a hyperlink, like its sibling cells) --
$newUrl = "https://pic1.imgdb.cn/item/694fe8e5161224305eb30b3e.jpg"
$ws.Range("C20").Value = $newUrl
$ws.Hyperlinks.Add($ws.Range("C20"), $newUrl)

# --- Restore the selection state left by the editing session ---------------
$ws.Range("C55").Select()
